$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("Main")
$fm = $wb.Worksheets.Item("Financial Model")

$fm.Range("X21").Value = 0.11
$fm.Range("X7").Formula = "=X5*0.62"
$fm.Range("Y7").Formula = "=Y5*0.63"
$fm.Range("Z7:AM7").Formula = "=Z5*0.63"

$main.Range("C6").Value = 2.4350000000000001

Write-Host "done"
